$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 186, shifting existing rows 186:225 down to 187:226
$ws.Rows.Item(186).Insert()

# Populate the new row 186 with the new data record
$ws.Range("A186").Value = 4
$ws.Range("B186").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C186").Value = "Los Lagos"
$ws.Range("D186").Value = 44782
$ws.Range("E186").Value = 10
$ws.Range("F186").Value = "Fruta"
$ws.Range("G186").Value = 100108
$ws.Range("H186").Value = "Tropicales y subtropicales"
$ws.Range("I186").Value = 100108002
$ws.Range("J186").Value = "Mango"
$ws.Range("K186").Value = "Sin especificar"
$ws.Range("L186").Value = "Primera"
$ws.Range("M186").Value = 160
$ws.Range("N186").Value = 13000
$ws.Range("O186").Value = 14000
$ws.Range("P186").Value = 13500
$ws.Range("Q186").Value = "$/bandeja 4 kilos"
$ws.Range("R186").Value = "Brasil"
$ws.Range("S186").Value = 3375
$ws.Range("T186").Value = 4
